$wb = $excel.ActiveWorkbook

# --- Sheet "Cadastro": update the login value in A2 ---
$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsCadastro.Range("A2").Value = "caiquasantanadeoliveira"

# --- Sheet "Pesquisa pagina inicial": new product report rows ---
$wsPesquisa = $wb.Worksheets.Item("Pesquisa pagina inicial")

# Row 2: new product (HP ZBook ...), rendered like a pasted web-report
# entry -> Consolas 9pt, color #222222, underlined.
$wsPesquisa.Range("A2").ClearFormats()
$wsPesquisa.Range("A2").Value = "HP ZBook 17 G2 Mobile Workstation"
$wsPesquisa.Range("A2").Font.Name = "Consolas"
$wsPesquisa.Range("A2").Font.Size = 9
$wsPesquisa.Range("A2").Font.Color = 2236962
$wsPesquisa.Range("A2").Font.Family = 3
$wsPesquisa.Range("A2").Font.Underline = $true

# Row 3: second product (Mouse optico ...) same report styling, no underline.
$wsPesquisa.Range("A3").ClearFormats()
$wsPesquisa.Range("A3").Value = "Mouse óptico USB com 3 botões HP"
$wsPesquisa.Range("A3").Font.Name = "Consolas"
$wsPesquisa.Range("A3").Font.Size = 9
$wsPesquisa.Range("A3").Font.Color = 2236962
$wsPesquisa.Range("A3").Font.Family = 3
